# Auto-generated script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.600.08"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").Value = "2.582.30"
$ws.Range("E3").Value = "  -2.70%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.580"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.77"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.22%  "
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("E11").Value = "  +3.59%  "
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("D13").Value = "3.036.65"
$ws.Range("E13").Value = "  -2.93%  "
$ws.Range("D14").Value = "58.500.05"
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("D16").Value = "2.579.36"
$ws.Range("E16").Value = "  -2.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000132"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.90%  "
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "334.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.05%  "
$ws.Range("E20").Value = "  -3.26%  "
$ws.Range("E21").Value = "  -4.19%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.423"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -5.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.11%  "
$ws.Range("D28").Value = "0.0₃0740"
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("E33").Value = "  -0.76%  "
$ws.Range("E34").Value = "  -3.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.850"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.823"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.64%  "
$ws.Range("E38").Value = "  -3.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "279.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.06%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  -2.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0941"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("E45").Value = "  -2.73%  "
$ws.Range("E46").Value = "  -5.43%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Value = "1.902.02"
$ws.Range("E48").Value = "  -4.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.08%  "
$ws.Range("E50").Value = "  -3.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "109.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.33%  "
